# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right after "2021-Q4" (before "总计"),
#    populated with the per-fund holding breakdown for that quarter.
# 2) Prepend a new row to the "总计" (totals) summary sheet with the
#    aggregated 2022-Q1 figures, shifting the older rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Header row
$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

$hdr = $newSheet.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Index column (A) formatting to match other quarter sheets
$idx = $newSheet.Range("A2:A6")
$idx.Font.Bold = $true
$idx.Borders.LineStyle = 1
$idx.HorizontalAlignment = -4108
$idx.VerticalAlignment = -4160

# Text columns B:G must stay text (fund codes / decimal-looking figures
# are stored as text in this workbook, not numbers)
$textCols = $newSheet.Range("B2:G6")
$textCols.NumberFormat = "@"

# Row 2 - 515210 国泰中证钢铁ETF
$newSheet.Range("A2").Value2 = 0
$newSheet.Range("B2").Value2 = "515210"
$newSheet.Range("C2").Value2 = "国泰中证钢铁ETF"
$newSheet.Range("D2").Value2 = "16.24"
$newSheet.Range("E2").Value2 = "99.25"
$newSheet.Range("F2").Value2 = "3.15"
$newSheet.Range("G2").Value2 = "0.5116"
$newSheet.Range("H2").Value2 = 8

# Row 3 - 502023 鹏华国证钢铁行业指数（LOF）
$newSheet.Range("A3").Value2 = 1
$newSheet.Range("B3").Value2 = "502023"
$newSheet.Range("C3").Value2 = "鹏华国证钢铁行业指数（LOF）"
$newSheet.Range("D3").Value2 = "15.55"
$newSheet.Range("E3").Value2 = "94.76"
$newSheet.Range("F3").Value2 = "2.76"
$newSheet.Range("G3").Value2 = "0.4292"
$newSheet.Range("H3").Value2 = 9

# Row 4 - 168203 中融国证钢铁行业指数
$newSheet.Range("A4").Value2 = 2
$newSheet.Range("B4").Value2 = "168203"
$newSheet.Range("C4").Value2 = "中融国证钢铁行业指数"
$newSheet.Range("D4").Value2 = "4.30"
$newSheet.Range("E4").Value2 = "92.58"
$newSheet.Range("F4").Value2 = "2.68"
$newSheet.Range("G4").Value2 = "0.1152"
$newSheet.Range("H4").Value2 = 9

# Row 5 - 013802 财通资管中证钢铁指数A
$newSheet.Range("A5").Value2 = 3
$newSheet.Range("B5").Value2 = "013802"
$newSheet.Range("C5").Value2 = "财通资管中证钢铁指数A"
$newSheet.Range("D5").Value2 = "0.11"
$newSheet.Range("E5").Value2 = "90.83"
$newSheet.Range("F5").Value2 = "2.85"
$newSheet.Range("G5").Value2 = "0.0031"
$newSheet.Range("H5").Value2 = 8

# Row 6 - 013803 财通资管中证钢铁指数C
$newSheet.Range("A6").Value2 = 4
$newSheet.Range("B6").Value2 = "013803"
$newSheet.Range("C6").Value2 = "财通资管中证钢铁指数C"
$newSheet.Range("D6").Value2 = "0.02"
$newSheet.Range("E6").Value2 = "90.83"
$newSheet.Range("F6").Value2 = "2.85"
$newSheet.Range("G6").Value2 = "0.0006"
$newSheet.Range("H6").Value2 = 8

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row on the "总计" sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Rows.Item(2).ClearFormats()

$totalSheet.Range("A2").Value2 = 0
$totalSheet.Range("B2").Value2 = "2022-Q1"
$totalSheet.Range("C2").Value2 = 5
$totalSheet.Range("D2").Value2 = 1.06

# Re-apply the index-column style (copied from the row right below, which
# still carries the original formatting) and renumber the shifted rows
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A3").Value2 = 1
$totalSheet.Range("A4").Value2 = 2
$totalSheet.Range("A5").Value2 = 3
$totalSheet.Range("A6").Value2 = 4

# Restore the originally active sheet (tab selection is not part of this change)
$wb.Worksheets.Item("2020-Q4").Activate()
